# Auto-generated edit script applying numeric corrections to Leve profit calculations
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 435.07693
$ws.Range("I28").Value = 333.45456
$ws.Range("K28").Value = 333.45456
$ws.Range("M28").Value = 151.54544
$ws.Range("H116").Value = 14399.6
$ws.Range("J116").Value = 17249.75
$ws.Range("L116").Value = 17249.75
$ws.Range("N116").Value = -24133.75
$ws.Range("H132").Value = 9844.923000000001
$ws.Range("I132").Value = 10498.167
$ws.Range("K132").Value = 31494.501
$ws.Range("M132").Value = -28964.501
$ws.Range("H138").Value = 2792.383
$ws.Range("J138").Value = 4478.4614
$ws.Range("L138").Value = 13435.3842
$ws.Range("N138").Value = -23715.3842

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1498
$ws.Range("I63").Value = 1498
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1498
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -812
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1498
$ws.Range("I66").Value = 1498
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 7490
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -4058
$ws.Range("N66").ClearContents()
$ws.Range("H110").Value = 725.43475
$ws.Range("I110").Value = 597.375
$ws.Range("J110").Value = 1018.1429
$ws.Range("K110").Value = 597.375
$ws.Range("L110").Value = 1018.1429
$ws.Range("M110").Value = 1447.625
$ws.Range("N110").Value = -5108.1429
$ws.Range("H122").Value = 2432.3914
$ws.Range("I122").Value = 2502.4
$ws.Range("J122").Value = 1965.6666
$ws.Range("K122").Value = 7507.200000000001
$ws.Range("L122").Value = 5896.9998
$ws.Range("M122").Value = -5057.200000000001
$ws.Range("N122").Value = -10796.9998
$ws.Range("H132").Value = 5616.5835
$ws.Range("I132").Value = 10099.75
$ws.Range("J132").Value = 3375
$ws.Range("K132").Value = 30299.25
$ws.Range("L132").Value = 10125
$ws.Range("M132").Value = -27769.25
$ws.Range("N132").Value = -15185

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33339396
$ws.Range("J20").Value = 2375.75
$ws.Range("L20").Value = 2375.75
$ws.Range("N20").Value = -2869.75
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H86").Value = 5559
$ws.Range("I86").Value = 5540.25
$ws.Range("J86").Value = 5589
$ws.Range("K86").Value = 5540.25
$ws.Range("L86").Value = 5589
$ws.Range("M86").Value = -4417.25
$ws.Range("N86").Value = -7835
$ws.Range("H89").Value = 5559
$ws.Range("I89").Value = 5540.25
$ws.Range("J89").Value = 5589
$ws.Range("K89").Value = 27701.25
$ws.Range("L89").Value = 27945
$ws.Range("M89").Value = -22085.25
$ws.Range("N89").Value = -39177

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 676.1818
$ws.Range("J94").Value = 749.4286
$ws.Range("L94").Value = 749.4286
$ws.Range("N94").Value = -1651.4286
$ws.Range("H132").Value = 3988.8438
$ws.Range("I132").Value = 3275.7827
$ws.Range("K132").Value = 9827.348100000001
$ws.Range("M132").Value = -7297.348100000001
$ws.Range("H134").Value = 3984.2415
$ws.Range("I134").Value = 4107.76
$ws.Range("K134").Value = 12323.28
$ws.Range("M134").Value = -9788.280000000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 4985
$ws.Range("J22").Value = 7998.75
$ws.Range("L22").Value = 23996.25
$ws.Range("N22").Value = -24334.25
$ws.Range("H27").Value = 4985
$ws.Range("J27").Value = 7998.75
$ws.Range("L27").Value = 23996.25
$ws.Range("N27").Value = -24200.25
$ws.Range("H32").Value = 9997.666999999999
$ws.Range("J32").Value = 9997.666999999999
$ws.Range("L32").Value = 29993.001
$ws.Range("N32").Value = -30559.001
$ws.Range("H46").Value = 604.875
$ws.Range("H81").Value = 2948.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2948.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 8845.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -11091.5
$ws.Range("H84").Value = 2948.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2948.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 26536.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -37768.5
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 3000
$ws.Range("M93").Value = -1128

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2163.1365
$ws.Range("I113").Value = 2164.182
$ws.Range("J113").Value = 2162.0908
$ws.Range("K113").Value = 2164.182
$ws.Range("L113").Value = 2162.0908
$ws.Range("M113").Value = 5.818000000000211
$ws.Range("N113").Value = -6502.0908

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2435.3333
$ws.Range("I7").Value = 1903.1666
$ws.Range("J7").Value = 3499.6667
$ws.Range("K7").Value = 1903.1666
$ws.Range("L7").Value = 3499.6667
$ws.Range("M7").Value = -1791.1666
$ws.Range("N7").Value = -3723.6667
$ws.Range("H46").Value = 2228.6155
$ws.Range("I46").Value = 1985.7142
$ws.Range("J46").Value = 2512
$ws.Range("K46").Value = 1985.7142
$ws.Range("L46").Value = 2512
$ws.Range("M46").Value = -1797.7142
$ws.Range("N46").Value = -2888
$ws.Range("H100").Value = 2807.25
$ws.Range("I100").Value = 3024.75
$ws.Range("J100").Value = 2372.25
$ws.Range("K100").Value = 3024.75
$ws.Range("L100").Value = 2372.25
$ws.Range("M100").Value = -2483.75
$ws.Range("N100").Value = -3454.25
$ws.Range("H126").Value = 2435.3333
$ws.Range("I126").Value = 1903.1666
$ws.Range("J126").Value = 3499.6667
$ws.Range("K126").Value = 5709.4998
$ws.Range("L126").Value = 10499.0001
$ws.Range("M126").Value = -3239.4998
$ws.Range("N126").Value = -15439.0001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4062
$ws.Range("I62").Value = 3333
$ws.Range("K62").Value = 3333
$ws.Range("M62").Value = -2709
$ws.Range("H65").Value = 4062
$ws.Range("I65").Value = 3333
$ws.Range("K65").Value = 16665
$ws.Range("M65").Value = -13545
$ws.Range("H74").Value = 5999.6665
$ws.Range("J74").Value = 5999.6665
$ws.Range("L74").Value = 5999.6665
$ws.Range("N74").Value = -7871.6665
$ws.Range("H77").Value = 5999.6665
$ws.Range("J77").Value = 5999.6665
$ws.Range("L77").Value = 17998.9995
$ws.Range("N77").Value = -27358.9995
$ws.Range("H124").Value = 500214.5
$ws.Range("J124").Value = 500214.5
$ws.Range("L124").Value = 500214.5
$ws.Range("N124").Value = -510034.5
$ws.Range("H132").Value = 3055.6775
$ws.Range("I132").Value = 3108.1724
$ws.Range("K132").Value = 9324.5172
$ws.Range("M132").Value = -6794.5172
